# Fixing conflicts in main for neutral-axis branch.
#
# The neutral-axis formula derivations for element #2 (the LE spar, row 3)
# and elements #3/#4 (the stringers, rows 4/5) were reworked to drop a
# small correction term (the "+(x/180)" / "+(x/90)" fudge factors) and to
# use the correct parallel-axis-theorem style rotated moment-of-inertia
# formulas for H3/I3/J3. Everything downstream (K:Y, AH:AM, AT, row 16
# totals, etc.) recalculates automatically from these root formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Root formula corrections (row 3: LE Spar) ---
$ws.Range("C3").Formula = "=SQRT(1+3.75^2)*0.0625"
$ws.Range("E3").Formula = "=(-SQRT(1+3.75^2)/2)*COS(ATAN(1/3.75))"
$ws.Range("H3").Formula = "=(((1/16)*(SQRT(1+3.75^2)))/24)*((1/16)^2+(SQRT(1+3.75^2))^2)+(((1/16)*(SQRT(1+3.75^2)))/24)*((1/16)^2-(SQRT(1+3.75^2))^2)*COS(2*ATAN(1/3.75))"
$ws.Range("I3").Formula = "=(((1/16)*(SQRT(1+3.75^2)))/24)*((1/16)^2+(SQRT(1+3.75^2))^2)-(((1/16)*(SQRT(1+3.75^2)))/24)*((1/16)^2-(SQRT(1+3.75^2))^2)*COS(2*ATAN(1/3.75))"
$ws.Range("J3").Formula = "=(((1/16)*(SQRT(1+3.75^2)))/24)*((1/16)^2-(SQRT(1+3.75^2))^2)*SIN(2*ATAN(1/3.75))"

# --- Root formula corrections (row 4: Top Stringer1) ---
$ws.Range("C4").Formula = "=1.5*0.125"
$ws.Range("D4").Formula = "=(-1.5)/2"

# --- Root formula corrections (row 5: Top Stringer2) ---
$ws.Range("C5").Formula = "=0.5*0.125"

# --- New note, merged AG8:AL8, centered ---
$ws.Range("AG8:AL8").HorizontalAlignment = -4108
$ws.Range("AG8").Value = "         Below values are assumed to be in the blue axis (auxillary axes) coord system below:"
$ws.Range("AG8:AL8").Merge()

# --- Selection moved to AJ7 ---
$ws.Range("AJ7").Select()
